# Updates the cryptos list prices (column D) and volume/percent-change
# values (column E) to the latest scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> [new Price (D) or $null if unchanged, new Volume(1h) (E)]
$changes = @{
    2  = @("62.812.49",       "  +2.42%  ")
    3  = @("2.945.99",        "  +0.50%  ")
    4  = @($null,             "  +0.00%  ")
    5  = @("592.18",          "  -0.44%  ")
    6  = @("147.32",          "  +2.66%  ")
    7  = @($null,             "  -0.06%  ")
    8  = @("2.945.04",        "  +0.55%  ")
    9  = @($null,             "  +0.83%  ")
    10 = @("7.03",            "  +1.21%  ")
    11 = @($null,             "  +5.25%  ")
    12 = @($null,             "  +0.11%  ")
    13 = @("0.0000233",       "  +4.21%  ")
    14 = @($null,             "  -2.34%  ")
    15 = @($null,             "  -0.92%  ")
    16 = @("3.432.95",        "  +0.46%  ")
    17 = @("62.795.64",       "  +2.39%  ")
    18 = @("6.66",            "  +0.40%  ")
    19 = @("2.938.24",        "  +0.31%  ")
    20 = @("438.52",          "  +1.21%  ")
    21 = @("13.40",           "  -1.02%  ")
    22 = @("0.664",           "  -1.19%  ")
    23 = @($null,             "  -0.86%  ")
    24 = @($null,             "  +2.97%  ")
    25 = @("80.72",           "  -0.91%  ")
    26 = @($null,             "  +0.81%  ")
    27 = @($null,             "  -2.00%  ")
    28 = @($null,             "  -0.03%  ")
    29 = @("2.22",            "  +1.00%  ")
    30 = @("7.27",            "  +5.69%  ")
    31 = @($null,             "  +0.37%  ")
    32 = @("0.0₃0983",       "  +12.69%  ")
    33 = @("26.34",           "  -1.12%  ")
    34 = @($null,             "  -0.41%  ")
    35 = @($null,             "  -0.06%  ")
    36 = @("0.990",           "  -2.09%  ")
    37 = @($null,             "  -0.30%  ")
    38 = @("3.02",            "  +2.04%  ")
    39 = @("49.59",           "  -0.30%  ")
    40 = @($null,             "  +1.33%  ")
    41 = @($null,             "  -3.79%  ")
    42 = @($null,             "  -0.79%  ")
    43 = @($null,             "  +0.22%  ")
    44 = @("39.23",           "  -6.89%  ")
    45 = @("2.701.86",        "  -0.05%  ")
    46 = @("135.13",          "  +1.26%  ")
    47 = @("0.0337",          "  -1.85%  ")
    48 = @("356.25",          "  -1.99%  ")
    50 = @($null,             "  -0.55%  ")
    51 = @("22.69",           "  -3.48%  ")
}

foreach ($row in $changes.Keys) {
    $pair = $changes[$row]
    $newPrice = $pair[0]
    $newVolume = $pair[1]

    if ($null -ne $newPrice) {
        $ws.Cells.Item($row, 4).Value2 = $newPrice
    }
    $ws.Cells.Item($row, 5).Value2 = $newVolume
}
